$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- column widths for the date/time columns (B, C, H) ---
$ws.Range("B1").ColumnWidth = 19.8333333333333
$ws.Range("C1").ColumnWidth = 19.8333333333333
$ws.Range("H1").ColumnWidth = 19.8333333333333

$dateFormat = "yyyy-mm-dd HH:mm:ss UTC"

# --- row 2 ---
$ws.Range("A2").Value = 1

$ws.Range("B2").Value = 44523.83829861111
$ws.Range("B2").NumberFormat = $dateFormat

$ws.Range("C2").Value = 44523.84128472222
$ws.Range("C2").NumberFormat = $dateFormat

$ws.Range("D2").Value = "IP Address"
$ws.Range("E2").Value = 100
$ws.Range("F2").Value = 258
$ws.Range("G2").Value = $true

$ws.Range("H2").Value = 44523.8412962963
$ws.Range("H2").NumberFormat = $dateFormat

$ws.Range("I2").Value = "1bimil"

$j2text = @"
ebola %>% 
  pivot_longer(Cases_Guinea:last_col()) %>% 
  separate(name, into = c("case_death", "country"), sep = "_") %>% 
  drop_na() %>% 
  pivot_wider(names_from = case_death, values_from = value)

"@
$ws.Range("J2").Value = $j2text

# the multi-line cell above makes the engine auto-grow row 2's height;
# put it back to the sheet's default so the row doesn't carry a custom height
$ws.Rows.Item(2).RowHeight = 15

# --- row 3 ---
$ws.Range("A3").Value = 2

$ws.Range("B3").Value = 44523.10165509259
$ws.Range("B3").NumberFormat = $dateFormat

$ws.Range("C3").Value = 44523.10333333333
$ws.Range("C3").NumberFormat = $dateFormat

$ws.Range("D3").Value = "Spam"
$ws.Range("E3").Value = 50
$ws.Range("F3").Value = 145
$ws.Range("G3").Value = $false

$ws.Range("H3").Value = 44523.87923611111
$ws.Range("H3").NumberFormat = $dateFormat

$ws.Range("I3").Value = "2nesch"
